$wb = $excel.ActiveWorkbook

# --- 1. Rename the "root_hospital_*" sheets to their short names ---
$wb.Worksheets.Item("root_hospital_cities").Name = "cities"
$wb.Worksheets.Item("root_hospital_patients").Name = "patients"
$wb.Worksheets.Item("root_hospital_users").Name = "users"

$wsCities = $wb.Worksheets.Item("cities")
$wsPatients = $wb.Worksheets.Item("patients")
$wsUsers = $wb.Worksheets.Item("users")
$wsAttr = $wb.Worksheets.Item("attributes")
$wsEnt = $wb.Worksheets.Item("entities")

# --- 2. Update "attributes" sheet data: entity / refEntity columns that used
#        to reference the long "root_hospital_*" sheet names now reference
#        the new short names ---
$wsAttr.Range("B2").Value = "cities"
$wsAttr.Range("B3").Value = "cities"
$wsAttr.Range("B4").Value = "cities"
$wsAttr.Range("B8").Value = "patients"
$wsAttr.Range("B9").Value = "patients"
$wsAttr.Range("E9").Value = "cities"
$wsAttr.Range("B10").Value = "patients"
$wsAttr.Range("E10").Value = "patients"
$wsAttr.Range("B11").Value = "patients"
$wsAttr.Range("B12").Value = "users"
$wsAttr.Range("B13").Value = "users"

# --- 3. "entities" sheet: drop the "package" column (old column B) -- the
#        "packages" sheet that it used to reference is being removed, so
#        "extends" / "abstract" / "description" shift one column left ---
[void]$wsEnt.Columns.Item(2).Delete()

# --- 4. Remove the now-unused "packages" sheet entirely ---
[void]$wb.Worksheets.Item("packages").Delete()

# --- 5. Update sheet selections to match the saved view state ---
[void]$wsCities.Range("C21").Select()
[void]$wsPatients.Range("D3").Select()
[void]$wsUsers.Range("E2").Select()
[void]$wsEnt.Range("D10").Select()

# "attributes" becomes the active/selected sheet
[void]$wsAttr.Activate()
[void]$wsAttr.Range("E9").Select()
